# Insert a new data row at sheet row 80 (pushing existing rows 80-116 down
# to 81-117) and populate it with the new "Perfection" / "Provincia de
# Huasco" record dated 44719, per the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(80).Insert()

$ws.Cells(80,1).Value = 4
$ws.Cells(80,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells(80,3).Value = "Los Lagos"
$ws.Cells(80,4).Value = 44719
$ws.Cells(80,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells(80,5).Value = 10
$ws.Cells(80,6).Value = 100112022
$ws.Cells(80,7).Value = "Arveja Verde"
$ws.Cells(80,8).Value = "Perfection"
$ws.Cells(80,9).Value = "Primera"
$ws.Cells(80,10).Value = 35
$ws.Cells(80,11).Value = 47000
$ws.Cells(80,12).Value = 47000
$ws.Cells(80,13).Value = 47000
$ws.Cells(80,14).Value = "$/malla 25 kilos"
$ws.Cells(80,15).Value = "Provincia de Huasco"
$ws.Cells(80,16).Value = 1880
$ws.Cells(80,17).Value = 25
$ws.Cells(80,18).Value = "Hortaliza"
